# CIERRE 5 abr 22
# - Switch the active/selected tab from "ARQUITECTO" to "VALES DE INSENTIVOS"
# - Update the incentive-month text from FEBRERO to MARZO on the
#   "VALES DE INSENTIVOS" sheet
# - Let the TODAY() volatile formulas on both vale sheets recompute to the
#   current (pinned) date

$wb = $excel.ActiveWorkbook

$wsVales = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# Update the "PAGO DE INCENTIVO DEL MES DE ..." text on the VALES sheet.
$wsVales.Range("A4").Value = "PAGO DE INCENTIVO DEL MES DE MARZO  2022"

# Recalculate so the TODAY() formulas (A11 on both sheets) pick up the
# current date.
$wb.Application.Calculate()

# Make "VALES DE INSENTIVOS" the active/selected sheet (was "ARQUITECTO").
$wsVales.Activate()
